# #35257 modify budzik boards
# Adds a new "duck size" board row (mała kaczka / duża kaczka) below the
# existing shape rows, with column C repeating the "correct" answer (A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "mała kaczka"
$ws.Range("B4").Value = "duża kaczka"
$ws.Range("C4").Value = "mała kaczka"
